$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the former row 141, pushing
# every subsequent record (old rows 141-209) down by one row.
$ws.Rows.Item(141).Insert()

$ws.Cells.Item(141, 1).Value = 3
$ws.Cells.Item(141, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(141, 3).Value = "Coquimbo"
$ws.Cells.Item(141, 4).Value = 44839
$ws.Cells.Item(141, 5).Value = 5
$ws.Cells.Item(141, 6).Value = 100112026
$ws.Cells.Item(141, 7).Value = "Haba"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 85
$ws.Cells.Item(141, 11).Value = 9000
$ws.Cells.Item(141, 12).Value = 10000
$ws.Cells.Item(141, 13).Value = 9529
$ws.Cells.Item(141, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(141, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(141, 16).Value = 381
$ws.Cells.Item(141, 17).Value = 25
$ws.Cells.Item(141, 18).Value = "Hortaliza"
